$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# A7: numeric id update
$ws.Range("A7").Value = 112222968

# B7: numeric id update
$ws.Range("B7").Value = 57620

# I7: text cell holding a numeric-looking value ("1" -> "2"); force text
# storage (the source cell is a text cell, not a number) then restore the
# default "Normal" style so no stray formatting is left behind.
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "2"
$ws.Range("I7").Style = "Normal"

# L7: text update
$ws.Range("L7").Value = "hane"
